# Apply the Sat May  6 16:11:50 UTC 2023 "Updated cryptos list" refresh to Sheet1.
# Each data row (2-51) gets new Price (D) / Volume(1h) (E) readings; rows 43-44
# additionally swap which coin (Bitcoin... PEPE/Algorand) occupies that rank, so
# Coin (B) and Link (C) are rewritten there too.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '28.666.17'
$ws.Range("E2").Value = '  -2.92%  '

# Row 3
$ws.Range("D3").Value = '1.882.79'
$ws.Range("E3").Value = '  -5.27%  '

# Row 4
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.37%  '

# Row 5
$ws.Range("D5").Value = '321.08'
$ws.Range("E5").Value = '  -2.22%  '

# Row 6
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.27%  '

# Row 7
$ws.Range("E7").Value = '  -1.92%  '

# Row 8
$ws.Range("D8").Value = '0.3772'
$ws.Range("E8").Value = '  -4.56%  '

# Row 9
$ws.Range("D9").Value = '45.38'
$ws.Range("E9").Value = '  -1.56%  '

# Row 10
$ws.Range("D10").Value = "'0.07680"
$ws.Range("E10").Value = '  -2.82%  '

# Row 11
$ws.Range("D11").Value = '0.9562'
$ws.Range("E11").Value = '  -4.90%  '

# Row 12
$ws.Range("E12").Value = '  -2.54%  '

# Row 13
$ws.Range("D13").Value = '1.891.42'
$ws.Range("E13").Value = '  -4.50%  '

# Row 14
$ws.Range("D14").Value = "'6.920"
$ws.Range("E14").Value = '  -4.15%  '

# Row 15
$ws.Range("D15").Value = '5.636'
$ws.Range("E15").Value = '  -3.87%  '

# Row 16
$ws.Range("D16").Value = '0.07019'
$ws.Range("E16").Value = '  -1.59%  '

# Row 17
$ws.Range("D17").Value = '1.003'
$ws.Range("E17").Value = '  +0.23%  '

# Row 18
$ws.Range("D18").Value = '82.53'
$ws.Range("E18").Value = '  -7.02%  '

# Row 19
$ws.Range("D19").Value = '0.000009471'
$ws.Range("E19").Value = '  -5.20%  '

# Row 20
$ws.Range("D20").Value = '16.56'
$ws.Range("E20").Value = '  -3.48%  '

# Row 21
$ws.Range("E21").Value = '  -0.77%  '

# Row 22
$ws.Range("D22").Value = '28.624.34'
$ws.Range("E22").Value = '  -3.35%  '

# Row 23
$ws.Range("D23").Value = '5.284'
$ws.Range("E23").Value = '  -4.84%  '

# Row 24
$ws.Range("D24").Value = '10.84'
$ws.Range("E24").Value = '  -3.61%  '

# Row 25
$ws.Range("D25").Value = '2.116.93'
$ws.Range("E25").Value = '  -4.11%  '

# Row 26
$ws.Range("D26").Value = '2.082'
$ws.Range("E26").Value = '  -1.83%  '

# Row 27
$ws.Range("D27").Value = '154.64'
$ws.Range("E27").Value = '  -2.07%  '

# Row 28
$ws.Range("D28").Value = '18.83'
$ws.Range("E28").Value = '  -3.93%  '

# Row 29
$ws.Range("D29").Value = "'5.610"
$ws.Range("E29").Value = '  -6.66%  '

# Row 30
$ws.Range("D30").Value = '116.45'
$ws.Range("E30").Value = '  -3.07%  '

# Row 31
$ws.Range("D31").Value = '1.799'
$ws.Range("E31").Value = '  -5.51%  '

# Row 32
$ws.Range("D32").Value = "'0.09200"
$ws.Range("E32").Value = '  -2.19%  '

# Row 33
$ws.Range("D33").Value = '0.8408'
$ws.Range("E33").Value = '  -6.63%  '

# Row 34
$ws.Range("D34").Value = '5.038'
$ws.Range("E34").Value = '  -4.44%  '

# Row 35
$ws.Range("D35").Value = "'1.240"
$ws.Range("E35").Value = '  -8.56%  '

# Row 36
$ws.Range("D36").Value = '2.929'
$ws.Range("E36").Value = '  -7.66%  '

# Row 37
$ws.Range("E37").Value = '  -3.33%  '

# Row 38
$ws.Range("E38").Value = '  -3.04%  '

# Row 39
$ws.Range("D39").Value = '1.001'
$ws.Range("E39").Value = '  +0.30%  '

# Row 40
$ws.Range("D40").Value = '0.02012'
$ws.Range("E40").Value = '  -5.55%  '

# Row 41
$ws.Range("D41").Value = '7.406'
$ws.Range("E41").Value = '  -5.70%  '

# Row 42
$ws.Range("D42").Value = '0.5449'
$ws.Range("E42").Value = '  -5.73%  '

# Row 43
$ws.Range("B43").Value = 'PEPE'
$ws.Range("C43").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D43").Value = '0.000003012'
$ws.Range("E43").Value = '  -20.36%  '

# Row 44
$ws.Range("B44").Value = 'Algorand'
$ws.Range("C44").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D44").Value = '0.1739'
$ws.Range("E44").Value = '  -4.75%  '

# Row 45
$ws.Range("D45").Value = "'9.130"
$ws.Range("E45").Value = '  -7.29%  '

# Row 46
$ws.Range("D46").Value = '2.676'
$ws.Range("E46").Value = '  +1.38%  '

# Row 47
$ws.Range("D47").Value = '0.5135'
$ws.Range("E47").Value = '  -4.70%  '

# Row 48
$ws.Range("D48").Value = '11.12'
$ws.Range("E48").Value = '  -8.53%  '

# Row 49
$ws.Range("D49").Value = '2.067'
$ws.Range("E49").Value = '  -7.43%  '

# Row 50
$ws.Range("D50").Value = '0.06748'
$ws.Range("E50").Value = '  -3.49%  '

# Row 51
$ws.Range("D51").Value = '110.88'
$ws.Range("E51").Value = '  -3.33%  '
